# Fixed update to excel issue
# - Rename "Requested quantity" headers to Weekly_PO_Qty / Monthly_PO_Qty
# - Add a new "PO Forecast" sheet with forecasted PO quantities (Prophet-style output)

$wb = $excel.ActiveWorkbook

# --- 1. Rename header on "Weekly Quantity" sheet (sheet1, B1) ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Rename header on "Monthly Trend" sheet (sheet2, B1) ---
$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the same header formatting already used on the other sheets (style index 1),
# instead of letting Excel mint a brand-new style entry.
$wsWeekly.Range("A1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null

# Forecast data rows 2-22: ds, PO_Forecast, yhat_lower, yhat_upper
$poData = @(
    @(44934.99999999999, 17, 8.675326675590361, 24.2884299137734),
    @(44941.99999999999, 16, 8.6162336840191, 23.96371238995656),
    @(44948.99999999999, 16, 8.354849378888785, 23.3886510284048),
    @(44955.99999999999, 15, 7.75437751811768, 22.33435591360416),
    @(44962.99999999999, 15, 7.253650646242308, 22.27138848012639),
    @(44976.99999999999, 14, 6.201342033874662, 20.42155711243951),
    @(44983.99999999999, 13, 5.396027936575519, 20.68104937471269),
    @(44990.99999999999, 13, 5.281057636172461, 20.00770483354632),
    @(44997.99999999999, 12, 4.569382495350565, 18.96802067851928),
    @(45004.99999999999, 12, 4.668241743309831, 18.78935623292343),
    @(45018.99999999999, 10, 3.466791282148916, 17.71560426357298),
    @(45039.99999999999, 9, 1.895667895547765, 16.220501842119),
    @(45074.99999999999, 6, -1.038365487615617, 13.79161933442399),
    @(45081.99999999999, 6, -1.966562534849354, 13.28698215125991),
    @(45088.99999999999, 5, -2.369452416102429, 13.05531543245993),
    @(45095.99999999999, 5, -2.865403694385821, 11.91593779355884),
    @(45102.99999999999, 4, -2.918316650595261, 11.31390359331273),
    @(45109.99999999999, 4, -3.886363582542864, 11.17607677276602),
    @(45116.99999999999, 3, -4.339313485876685, 10.53246303253703),
    @(45123.99999999999, 3, -4.854899896194637, 9.752994026682368),
    @(45130.99999999999, 2, -4.834569675917209, 9.812775050451609)
)

$row = 2
foreach ($rec in $poData) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# Apply the same date/time formatting (style index 2) used for column A on the other sheets
$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A22").PasteSpecial(-4122) | Out-Null

$wsForecast.Range("A1").Select()
